$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Uzimtas?" task name to "Create UI for the program"
$ws.Range("C13").Value = "Create UI for the program"

# Clear the now-removed task row (row 14, columns A-D)
$ws.Range("A14:D14").ClearContents()

# Update the active selection
$ws.Range("E9").Select()

$wb.Save()
